# From v1.0.1 to v1.0.2
# Swap the "Step 2" content of TC2 and TC3:
#  - TC2's step 2 (row 20) previously described the search/filter action; it
#    must now describe the cancellation action.
#  - TC3's step 2 (row 28) previously described the cancellation action; it
#    must now describe the search/filter action.
# The "Test Case ID" labels (TC2 in B15, TC3 in B23) and all other rows stay
# untouched; only the text in B20/D20 and B28/D28 is exchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$searchStep  = "Chefe Indica alguns parâmetros específicos para a busca; Informa o nome do beneficiário; Filtra a listagem de solicitações."
$searchResult = "SYSTEM Exibe uma nova listagem de solicitações, de acordo com os filtros informados pelo usuário."

$cancelStep   = "Chefe Clica para realizar o cancelamento de uma diária."
$cancelResult = "SYSTEM Verifica que a solicitação está em situação SOLICITADA; Exibe mensagem de confirmação (MSG987 - Cancelar solicitação de diária) para o usuário (que deve confirmar); Cancela a diária, mudando sua situação para CANCELADA (ver diagrama de estados da diária)."

# TC2 block (rows 15-20) -> step 2 becomes the cancellation text
$ws.Range("B20").Value = $cancelStep
$ws.Range("D20").Value = $cancelResult

# TC3 block (rows 23-28) -> step 2 becomes the search/filter text
$ws.Range("B28").Value = $searchStep
$ws.Range("D28").Value = $searchResult
